# Apply the "Removed 1 deleted question" edit to the county_councils sheet.
#
# Summary of the change:
#  - One sub-question (column H / s4_coms) was removed from the scoring
#    rubric, so every council's s4_coms score is rescaled from an "out of
#    10" fraction to an "out of 9" fraction (weighted_total, column N,
#    is recomputed accordingly).
#  - Because the weighted_total changed, the row ordering (sorted
#    descending by weighted_total) shuffled for a couple of adjacent
#    councils that swapped rank:
#      * Surrey County Council  <-> North Yorkshire County Council
#      * East Sussex County Council <-> Derbyshire County Council
#    so those two row-pairs have their entire row contents swapped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2 (Somerset) ----
$ws.Range("H2").Value = 0.7777777777777778
$ws.Range("N2").Value = 0.6280952380952382

# ---- Row 3 (Cambridgeshire) ----
$ws.Range("H3").Value = 0.5555555555555556
$ws.Range("N3").Value = 0.5594047619047618

# ---- Row 4 (Suffolk) ----
$ws.Range("H4").Value = 0.6666666666666666
$ws.Range("N4").Value = 0.5328571428571429

# ---- Row 5 (Kent) ----
$ws.Range("H5").Value = 0.6666666666666666
$ws.Range("N5").Value = 0.5297619047619048

# ---- Row 6 (Hampshire) ----
$ws.Range("H6").Value = 0.6666666666666666
$ws.Range("N6").Value = 0.5289285714285714

# ---- Row 7 (Hertfordshire) ----
$ws.Range("H7").Value = 0.6666666666666666
$ws.Range("N7").Value = 0.4970238095238095

# ---- Row 8: was Surrey, becomes North Yorkshire ----
$ws.Range("C8").Value = "North Yorkshire County Council"
$ws.Range("D8").Value = "NYK"
$ws.Range("E8").Value = 0.5238095238095238
$ws.Range("F8").Value = 0.05555555555555555
$ws.Range("G8").Value = 0.8
$ws.Range("H8").Value = 0.6666666666666666
$ws.Range("I8").Value = 0.4
$ws.Range("J8").Value = 0.25
$ws.Range("K8").Value = 0.2
$ws.Range("L8").Value = 0.8
$ws.Range("M8").Value = 0.5
$ws.Range("N8").Value = 0.484404761904762
$ws.Range("O8").Value = 2
$ws.Range("P8").Value = "under 800k"
$ws.Range("Q8").Value = "Yorkshire and The Humber"

# ---- Row 9 (Leicestershire, stays put) ----
$ws.Range("H9").Value = 0.4444444444444444
$ws.Range("N9").Value = 0.483452380952381

# ---- Row 10: was North Yorkshire, becomes Surrey ----
$ws.Range("C10").Value = "Surrey County Council"
$ws.Range("D10").Value = "SRY"
$ws.Range("E10").Value = 0.3333333333333333
$ws.Range("F10").Value = 0.9444444444444444
$ws.Range("G10").Value = 0.4
$ws.Range("H10").Value = 0.4444444444444444
$ws.Range("I10").Value = 0.6
$ws.Range("J10").Value = 0.5
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0.6
$ws.Range("M10").Value = 0.25
$ws.Range("N10").Value = 0.4758333333333333
$ws.Range("O10").Value = 2
$ws.Range("P10").Value = "1m +"
$ws.Range("Q10").Value = "South East"

# ---- Row 11 (Devon) ----
$ws.Range("H11").Value = 0.8888888888888888
$ws.Range("N11").Value = 0.4679761904761905

# ---- Row 12 (Oxfordshire) ----
$ws.Range("H12").Value = 0.5555555555555556
$ws.Range("N12").Value = 0.4264285714285714

# ---- Row 13 (Worcestershire, stays put) ----
$ws.Range("H13").Value = 0.3333333333333333
$ws.Range("N13").Value = 0.4038095238095239

# ---- Row 14: was East Sussex, becomes Derbyshire ----
$ws.Range("C14").Value = "Derbyshire County Council"
$ws.Range("D14").Value = "DBY"
$ws.Range("E14").Value = 0.3809523809523809
$ws.Range("F14").Value = 0.2222222222222222
$ws.Range("G14").Value = 0.4
$ws.Range("H14").Value = 0.4444444444444444
$ws.Range("I14").Value = 0.6
$ws.Range("J14").Value = 0.5
$ws.Range("K14").Value = 0.2
$ws.Range("L14").Value = 0.6
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = 0.3821428571428572
$ws.Range("O14").Value = 2
$ws.Range("P14").Value = "800k - 1m"
$ws.Range("Q14").Value = "West Midlands"

# ---- Row 15: was Derbyshire, becomes East Sussex ----
$ws.Range("C15").Value = "East Sussex County Council"
$ws.Range("D15").Value = "ESX"
$ws.Range("E15").Value = 0.4285714285714285
$ws.Range("F15").Value = 0.1111111111111111
$ws.Range("G15").Value = 0.6
$ws.Range("H15").Value = 0.3333333333333333
$ws.Range("I15").Value = 0.8
$ws.Range("J15").Value = 0.75
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 0.75
$ws.Range("N15").Value = 0.3759523809523809
$ws.Range("O15").Value = 1
$ws.Range("P15").Value = "under 800k"
$ws.Range("Q15").Value = "South East"

# ---- Row 16 (West Sussex) ----
$ws.Range("H16").Value = 0.4444444444444444
$ws.Range("N16").Value = 0.3469047619047619

# ---- Row 17 (Gloucestershire) ----
$ws.Range("H17").Value = 0.5555555555555556
$ws.Range("N17").Value = 0.3288095238095238

# ---- Row 18 (Warwickshire) ----
$ws.Range("H18").Value = 0.6666666666666666
$ws.Range("N18").Value = 0.2508333333333334

# ---- Row 19 (Staffordshire) ----
$ws.Range("H19").Value = 0.4444444444444444
$ws.Range("N19").Value = 0.2445238095238095

# ---- Row 20 (Nottinghamshire) ----
$ws.Range("H20").Value = 0.2222222222222222
$ws.Range("N20").Value = 0.1952380952380953

# ---- Row 21 (Lincolnshire) ----
$ws.Range("H21").Value = 0.3333333333333333
$ws.Range("N21").Value = 0.1922619047619047

# ---- Row 22 (Cumbria) ----
$ws.Range("H22").Value = 0.2222222222222222
$ws.Range("N22").Value = 0.1554761904761905
